# Update column F (dSF) values to match repulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = -2
    3  = 1
    4  = 3
    5  = 2
    6  = -1
    7  = -1
    8  = -3
    9  = 1
    10 = 4
    11 = -4
    13 = -5
    14 = 1
    17 = -6
    18 = -3
    19 = -2
    20 = -3
    22 = 1
    23 = 2
    25 = -1
    26 = -6
    27 = -1
    30 = 3
    32 = 2
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 6).Value = $values[$row]
}
